$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2:G5").Value = 0.4821145882335259
$ws.Range("H2:H5").Value = 0.992
